$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet2" -> "payment-request"
$paymentSheet = $wb.Worksheets.Item("Sheet2")
$paymentSheet.Name = "payment-request"

# 2. Fix the stray curly-quote typo in the 'CASH' literal (row 3, col J).
#    Assigning a string that literally starts with an apostrophe through
#    .Value/.Value2 turns on the cell's quote-prefix (text) formatting,
#    which would needlessly change the cell's style index. Route the text
#    through a formula + paste-values so only the literal text changes.
$paymentSheet.Range("J3").Formula = "=""'CASH'"""
$paymentSheet.Range("J3").Copy()
$paymentSheet.Range("J3").PasteSpecial(-4163)

# 3. Add a new "type" row (row 4) describing each column's data type
$paymentSheet.Range("A4").Value = "Long"
$paymentSheet.Range("B4").Value = "String"
$paymentSheet.Range("C4").Value = "Date"
$paymentSheet.Range("D4").Value = "Date"
$paymentSheet.Range("E4").Value = "String"
$paymentSheet.Range("F4").Value = "String"
$paymentSheet.Range("G4").Value = "String"
$paymentSheet.Range("H4").Value = "String"
$paymentSheet.Range("I4").Value = "Long"
$paymentSheet.Range("J4").Value = "String"
$paymentSheet.Range("K4").Value = "Float"

# 4. Move the selection on the "payment-request" sheet to D28
[void]$paymentSheet.Range("D28").Select()
